{"js": "// The document has a stray paragraph right after the \"RUT\" (Ruth) book\n// Heading 2 that just contains an italic \"Ruthu\" run \u2014 a duplicate of the\n// proper \"Ruthu\" Heading 2 section further down. This paragraph needs to\n// be removed entirely (including its paragraph mark), leaving the \"RUT\"\n// heading followed directly by the next (whitespace-only) paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the Heading2 paragraph whose text is exactly \"RUT\".\nlet rutIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"RUT\") {\n    rutIndex = i;\n    break;\n  }\n}\n\nif (rutIndex === -1) {\n  throw new Error('Could not find the \"RUT\" heading paragraph.');\n}\n\n// The very next paragraph is the stray italic \"Ruthu\" paragraph that must\n// be deleted outright (paragraph mark included).\nconst strayParagraph = paragraphs.items[rutIndex + 1];\nstrayParagraph.load(\"text\");\nawait context.sync();\n\nif (strayParagraph.text !== \"Ruthu\") {\n  throw new Error(\n    'Unexpected paragraph after \"RUT\": \"' + strayParagraph.text + '\"'\n  );\n}\n\nstrayParagraph.delete();\nawait context.sync();\n", "ps1": "# The document has a stray paragraph right after the \"RUT\" (Ruth) book\n# Heading 2 that just contains an italic \"Ruthu\" run -- a duplicate of the\n# proper \"Ruthu\" Heading 2 section further down. This paragraph needs to\n# be removed entirely (including its paragraph mark), leaving the \"RUT\"\n# heading followed directly by the next (whitespace-only) paragraph.\n\n$d = $word.ActiveDocument\n\n# Find the Heading2 paragraph whose text is exactly \"RUT\".\n$rutIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($t -eq \"RUT\") {\n    $rutIndex = $i\n    break\n  }\n}\n\nif ($rutIndex -eq -1) {\n  throw \"Could not find the 'RUT' heading paragraph.\"\n}\n\n# The very next paragraph is the stray italic \"Ruthu\" paragraph that must\n# be deleted outright (paragraph mark included).\n$strayIndex = $rutIndex + 1\n$stray = $d.Paragraphs.Item($strayIndex)\n$strayText = $stray.Range.Text.TrimEnd([char]13, [char]7)\n\nif ($strayText -ne \"Ruthu\") {\n  throw \"Unexpected paragraph after 'RUT': '$strayText'\"\n}\n\n$stray.Range.Delete()\n"}
